# Daily attendance processing - 2026-01-03 11:53:30
# In the "Recorded By" column, swap the order of the two recorder names
# so that "dnasr281@gmail.com, System" becomes "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Text -eq $oldText) {
            $cell.Value = $newText
        }
    }
}
